# Exercice 2 : jour en français
# 1) Translate the English weekday names in "Données"!H2:H16 to French.
# 2) Add a leading "region" column to "Par région" (shifting quantite -> B).
# 3) Add a leading "produit" column to "Par produit" (shifting quantite -> B).

$wb = $excel.ActiveWorkbook

# --- 1) "Données" sheet: jour_semaine column (H) English -> French ---
$wsData = $wb.Worksheets.Item("Données")

$jours = @{
    "Monday"    = "Lundi"
    "Tuesday"   = "Mardi"
    "Wednesday" = "Mercredi"
    "Thursday"  = "Jeudi"
    "Friday"    = "Vendredi"
    "Saturday"  = "Samedi"
    "Sunday"    = "Dimanche"
}

for ($row = 2; $row -le 16; $row++) {
    $cell = $wsData.Cells.Item($row, 8)   # column H (jour_semaine)
    $en = $cell.Value()
    if ($jours.ContainsKey($en)) {
        $cell.Value = $jours[$en]
    }
}

# --- 2) "Par région" sheet: insert "region" column before "quantite" ---
$wsRegion = $wb.Worksheets.Item("Par région")

$regionValues = @("Est", "Non spécifié", "Nord", "Ouest", "Sud")

# Insert a new blank column at A; the former column A ("quantite") shifts to B.
$wsRegion.Columns.Item(1).Insert()

# Write the new header + region values into column A
$wsRegion.Cells.Item(1, 1).Value = "region"
for ($i = 0; $i -lt $regionValues.Length; $i++) {
    $wsRegion.Cells.Item($i + 2, 1).Value = $regionValues[$i]
}

# Give the whole new column the same look as the sheet's header style
# (bold, thin border all round, centered / top-aligned).
$rngA = $wsRegion.Range("A1:A6")
$rngA.Font.Bold = $true
$rngA.Borders.LineStyle = 1
$rngA.HorizontalAlignment = -4108   # xlCenter
$rngA.VerticalAlignment = -4160     # xlTop

# --- 3) "Par produit" sheet: insert "produit" column before "quantite" ---
$wsProduit = $wb.Worksheets.Item("Par produit")

$produitValues = @("Souris")

# Insert a new blank column at A; the former column A ("quantite") shifts to B.
$wsProduit.Columns.Item(1).Insert()

# Write the new header + produit values into column A
$wsProduit.Cells.Item(1, 1).Value = "produit"
for ($i = 0; $i -lt $produitValues.Length; $i++) {
    $wsProduit.Cells.Item($i + 2, 1).Value = $produitValues[$i]
}

# Give the whole new column the same look as the sheet's header style
# (bold, thin border all round, centered / top-aligned).
$rngA2 = $wsProduit.Range("A1:A2")
$rngA2.Font.Bold = $true
$rngA2.Borders.LineStyle = 1
$rngA2.HorizontalAlignment = -4108   # xlCenter
$rngA2.VerticalAlignment = -4160     # xlTop
